# Commit change using GitKraken
#
# The document ends with a single empty "ListParagraph" bullet (level 2,
# ilvl=1, numId=1). Replace it with four fully-formed bullets describing
# the Git Log command walkthrough and finishing with the GitKraken note.

$d = $word.ActiveDocument

# The target is the last paragraph in the document: an empty bullet at
# list level 2 (w:ilvl="1") that was left as a placeholder.
$target = $d.Paragraphs.Last

$newParagraphsXml = @'
<pkg:xmlData xmlns:pkg="http://schemas.openxmlformats.org/package/2006/content-types"><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Closer look at the Git Log Command</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Git log --abbrev-commit</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="1"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Git log – -</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>oneline</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:t>GitKraken</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> – GUI to commit using a graphic user interface.</w:t></w:r></w:p></pkg:xmlData>
'@

# Replacing the target range's XML with four <w:p> elements turns the one
# placeholder paragraph into the four bullets described above, in place.
$target.Range.InsertXML($newParagraphsXml)
